$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Sheet1 ("...Signal Detection_2024-09"): append summary row 102 with
#    column totals for G:J (hit / miss / fa / cr counts).
# ---------------------------------------------------------------------------
$ws1.Range("G102").Formula = "=SUM(G2:G101)"
$ws1.Range("H102").Formula = "=SUM(H2:H101)"
$ws1.Range("I102").Formula = "=SUM(I2:I101)"
$ws1.Range("J102").Formula = "=SUM(J2:J101)"

# ---------------------------------------------------------------------------
# 2. Add the new "Sheet1" worksheet right after the existing sheet - it will
#    become the active / selected tab, matching the workbook's new state.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

# Column widths (best-effort - engine rounds to its own pixel grid).
$ws2.Columns.Item(2).ColumnWidth = 22.8
$ws2.Columns.Item(3).ColumnWidth = 12.3

# --- Populate cells in the same order the strings were first authored so ---
# --- the rebuilt shared-string table lines up with the source workbook.  ---

# Row 4 / Row 10 header pair (top confusion-matrix + the raw counts table).
$ws2.Range("D4").Value = "respond no"
$ws2.Range("B5").Value = "signal present (hit)"
$ws2.Range("B6").Value = "signal present c"
$ws2.Range("D6").Value = "corect rejection"
$ws2.Range("C4").Value = "respond yes"
$ws2.Range("C6").Value = "FALSE alarm"
$ws2.Range("B11").Value = "signal presengt"
$ws2.Range("B12").Value = "signal absent"
$ws2.Range("H5").Value = "prop hit=hit/hit+miss"
$ws2.Range("H6").Value = "prop Fa=FA/FA+CORR REJECTION"
$ws2.Range("H8").Value = "D PRIME=z(prop hit)-z(prop fa)"
$ws2.Range("H9").Value = "c=-z(prop hit)+z(prop fa)/2"

# Re-used strings / labels already present in the shared-string table.
$ws2.Range("C5").Value = "hit"
$ws2.Range("D5").Value = "miss"
$ws2.Range("C10").Value = "respond yes"
$ws2.Range("D10").Value = "respond no"

# Raw counts (pulled from the totals row on the data sheet).
$ws2.Range("C11").Value = 40
$ws2.Range("D11").Value = 6
$ws2.Range("C12").Value = 33
$ws2.Range("D12").Value = 21

# Derived proportions / statistics.
$ws2.Range("M5").Formula = "=C11/(C11+D11)"
$ws2.Range("M6").Formula = "=C12/(C12+D12)"
$ws2.Range("M8").Formula = "=NORMSINV(M5)-NORMSINV(M6)"
$ws2.Range("M9").Formula = "=NORMSINV(M5)+NORMSINV(M6)/2"

# ---------------------------------------------------------------------------
# 3. View state: sheet1 scrolled near the bottom with K102 selected, new
#    Sheet1 active with M8 selected (matches the recorded workbook views).
# ---------------------------------------------------------------------------
$ws1.Range("K102").Select()
$ws2.Activate()
$ws2.Range("M8").Select()
